# "averaged ptm plot in nb"
# Build the percentage table (columns G:K) next to the existing counts
# table (A:E) on the "MF for python" sheet, replace the GO-id labels in
# column A with their human-readable names, and add a "total" row (11)
# with SUM() formulas.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("MF for python")
$ws.Activate()

# --- Column A: swap the GO-term id for its readable name -------------
$ws.Range("A2").Value = "ATP binding"
$ws.Range("A3").Value = "DNA binding"
$ws.Range("A4").Value = "metal ion binding"
$ws.Range("A5").Value = "oxidoreductase activity"
$ws.Range("A6").Value = "ATPase"
$ws.Range("A7").Value = "zinc ion binding"
$ws.Range("A8").Value = "RNA binding"
$ws.Range("A9").Value = "hydrolase activity"
$ws.Range("A10").Value = "transmembrane transporter activity"

# --- Row 1: mirror the header row (incl. formatting) into G1:K1 ------
$ws.Range("A1:E1").Copy()
[void]$ws.Range("G1").Select()
[void]$ws.Paste()
$excel.CutCopyMode = $false

# --- Column G: mirror the (new) names from column A -------------------
$ws.Range("G2").Value = $ws.Range("A2").Value()
$ws.Range("G3").Value = $ws.Range("A3").Value()
$ws.Range("G4").Value = $ws.Range("A4").Value()
$ws.Range("G5").Value = $ws.Range("A5").Value()
$ws.Range("G6").Value = $ws.Range("A6").Value()
$ws.Range("G7").Value = $ws.Range("A7").Value()
$ws.Range("G8").Value = $ws.Range("A8").Value()
$ws.Range("G9").Value = $ws.Range("A9").Value()
$ws.Range("G10").Value = $ws.Range("A10").Value()

# --- Row 11: totals per count column (B:E) -----------------------------
$ws.Range("A11").Value = "total"
$ws.Range("B11").Formula = "=SUM(B2:B10)"
$ws.Range("C11").Formula = "=SUM(C2:C10)"
$ws.Range("D11").Formula = "=SUM(D2:D10)"
$ws.Range("E11").Formula = "=SUM(E2:E10)"

# --- Columns H:K rows 2:10: percentage of each column's total ---------
for ($r = 2; $r -le 10; $r++) {
    $ws.Range("H$r").Formula = "=B$r/B`$11*100"
    $ws.Range("I$r").Formula = "=C$r/C`$11*100"
    $ws.Range("J$r").Formula = "=D$r/D`$11*100"
    $ws.Range("K$r").Formula = "=E$r/E`$11*100"
}

[void]$ws.Range("G1").Select()
